# [Kadastro App] Yeni kayit eklendi: 2957
# Adds a new record row (row 48) to both the "Kayitlar" master list sheet
# and the "Erdemli" district sheet, mirroring the existing data layout.

$wb = $excel.ActiveWorkbook

$newRow = @("2957", "2025-09-09", "Erdemli", "1", "ÇAP", "CEMAL TİMUROĞLU (K.Teknisyeni)")
$sheetNames = @("Kayitlar", "Erdemli")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    $rowRange = $ws.Range("A48:F48")

    # Force the cells to be stored as text (matching the rest of the
    # column, which stores every value - including numeric-looking ones
    # like the record id - as text) instead of letting Excel infer
    # numbers/dates from the values.
    $rowRange.NumberFormat = "@"

    $ws.Range("A48").Value = $newRow[0]
    $ws.Range("B48").Value = $newRow[1]
    $ws.Range("C48").Value = $newRow[2]
    $ws.Range("D48").Value = $newRow[3]
    $ws.Range("E48").Value = $newRow[4]
    $ws.Range("F48").Value = $newRow[5]

    # Restore the default "Normal" style so the new cells don't end up
    # with a distinct cell style compared to the rest of the sheet
    # (only the underlying stored type should differ).
    $rowRange.Style = "Normal"
}
